$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1376.4
$ws.Range("I38").Value = 129
$ws.Range("J38").Value = 2208
$ws.Range("K38").Value = 387
$ws.Range("L38").Value = 6624
$ws.Range("M38").Value = -15
$ws.Range("N38").Value = -7368

$ws.Range("H92").Value = 3996
$ws.Range("I92").Value = 4152
$ws.Range("J92").Value = 1500
$ws.Range("K92").Value = 4152
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = -2904
$ws.Range("N92").Value = -3996

$ws.Range("H101").Value = 1139.2727
$ws.Range("I101").Value = 1463.1428
$ws.Range("J101").Value = 572.5
$ws.Range("K101").Value = 4389.428400000001
$ws.Range("L101").Value = 1717.5
$ws.Range("M101").Value = -2767.428400000001
$ws.Range("N101").Value = -4961.5

$ws.Range("H116").Value = 3048
$ws.Range("I116").Value = 2719.56
$ws.Range("K116").Value = 2719.56
$ws.Range("M116").Value = 722.4400000000001

$ws.Range("H135").Value = 4439.971
$ws.Range("I135").Value = 2787.7727
$ws.Range("K135").Value = 25089.9543
$ws.Range("M135").Value = -22554.9543

$ws.Range("H138").Value = 3327.2354
$ws.Range("J138").Value = 2961.1082
$ws.Range("L138").Value = 8883.3246
$ws.Range("N138").Value = -19163.3246

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2885.7144
$ws.Range("I5").Value = 3333.3333
$ws.Range("K5").Value = 3333.3333
$ws.Range("M5").Value = -3221.3333

$ws.Range("H61").Value = 22366.666
$ws.Range("I61").Value = 13730.333
$ws.Range("K61").Value = 13730.333
$ws.Range("M61").Value = -13518.333

$ws.Range("H96").Value = 26335.875
$ws.Range("J96").Value = 26335.875
$ws.Range("L96").Value = 26335.875
$ws.Range("N96").Value = -31827.875

$ws.Range("H102").Value = 15083.0625
$ws.Range("I102").Value = 2755.3333
$ws.Range("J102").Value = 199999
$ws.Range("K102").Value = 2755.3333
$ws.Range("L102").Value = 199999
$ws.Range("M102").Value = -1133.3333
$ws.Range("N102").Value = -203243

$ws.Range("H122").Value = 2811.0212
$ws.Range("I122").Value = 1948.7742
$ws.Range("J122").Value = 4481.625
$ws.Range("K122").Value = 5846.3226
$ws.Range("L122").Value = 13444.875
$ws.Range("M122").Value = -3396.3226
$ws.Range("N122").Value = -18344.875

$ws.Range("H136").Value = 22366.666
$ws.Range("I136").Value = 13730.333
$ws.Range("K136").Value = 41190.999
$ws.Range("M136").Value = -38640.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2885.7144
$ws.Range("I4").Value = 3333.3333
$ws.Range("K4").Value = 3333.3333
$ws.Range("M4").Value = -3218.3333

$ws.Range("H86").Value = 4494.136
$ws.Range("I86").Value = 4734.273
$ws.Range("J86").Value = 4254
$ws.Range("K86").Value = 4734.273
$ws.Range("L86").Value = 4254
$ws.Range("M86").Value = -3611.273
$ws.Range("N86").Value = -6500

$ws.Range("H89").Value = 4494.136
$ws.Range("I89").Value = 4734.273
$ws.Range("J89").Value = 4254
$ws.Range("K89").Value = 23671.365
$ws.Range("L89").Value = 21270
$ws.Range("M89").Value = -18055.365
$ws.Range("N89").Value = -32502

$ws.Range("H99").Value = 722
$ws.Range("I99").Value = 718.5
$ws.Range("J99").Value = 750
$ws.Range("K99").Value = 718.5
$ws.Range("L99").Value = 750
$ws.Range("M99").Value = 779.5
$ws.Range("N99").Value = -3746

$ws.Range("H105").Value = 4082.875
$ws.Range("I105").Value = 4082.875
$ws.Range("K105").Value = 4082.875
$ws.Range("M105").Value = -2335.875

$ws.Range("H107").Value = 8489.5
$ws.Range("I107").Value = 8877.223
$ws.Range("K107").Value = 8877.223
$ws.Range("M107").Value = -6957.223

$ws.Range("H134").Value = 10451.667
$ws.Range("I134").Value = 2726.2
$ws.Range("K134").Value = 8178.599999999999
$ws.Range("M134").Value = -5643.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28459.45
$ws.Range("I31").Value = 16391.428
$ws.Range("K31").Value = 16391.428
$ws.Range("M31").Value = -16096.428

$ws.Range("H34").Value = 28459.45
$ws.Range("I34").Value = 16391.428
$ws.Range("K34").Value = 16391.428
$ws.Range("M34").Value = -16189.428

$ws.Range("H58").Value = 34998.418
$ws.Range("I58").Value = 22629.334
$ws.Range("K58").Value = 22629.334
$ws.Range("M58").Value = -22426.334

$ws.Range("H62").Value = 11273.866
$ws.Range("I62").Value = 11932.833
$ws.Range("J62").Value = 10834.556
$ws.Range("K62").Value = 11932.833
$ws.Range("L62").Value = 10834.556
$ws.Range("M62").Value = -11308.833
$ws.Range("N62").Value = -12082.556

$ws.Range("H65").Value = 11273.866
$ws.Range("I65").Value = 11932.833
$ws.Range("J65").Value = 10834.556
$ws.Range("K65").Value = 59664.165
$ws.Range("L65").Value = 54172.78
$ws.Range("M65").Value = -56544.165
$ws.Range("N65").Value = -60412.78

$ws.Range("H70").Value = 16666.666
$ws.Range("J70").Value = 16666.666
$ws.Range("L70").Value = 16666.666
$ws.Range("N70").Value = -17296.666

$ws.Range("H73").Value = 16666.666
$ws.Range("J73").Value = 16666.666
$ws.Range("L73").Value = 16666.666
$ws.Range("N73").Value = -18850.666

$ws.Range("H136").Value = 34998.418
$ws.Range("I136").Value = 22629.334
$ws.Range("K136").Value = 67888.00199999999
$ws.Range("M136").Value = -65338.00199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 5249.5
$ws.Range("I22").Value = 5249.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 15748.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -15579.5
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 5249.5
$ws.Range("I27").Value = 5249.5
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 15748.5
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -15646.5
$ws.Range("N27").ClearContents()

$ws.Range("H57").Value = 5000
$ws.Range("J57").Value = 5000
$ws.Range("L57").Value = 15000
$ws.Range("N57").Value = -16118

$ws.Range("H140").Value = 2789.1667
$ws.Range("I140").Value = 1352.8572
$ws.Range("K140").Value = 4058.5716
$ws.Range("M140").Value = 1121.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 34248.5
$ws.Range("J92").Value = 34248.5
$ws.Range("L92").Value = 34248.5
$ws.Range("N92").Value = -37992.5

$ws.Range("H102").Value = 6014.231
$ws.Range("I102").Value = 2548.25
$ws.Range("J102").Value = 11559.8
$ws.Range("K102").Value = 2548.25
$ws.Range("L102").Value = 11559.8
$ws.Range("M102").Value = -926.25
$ws.Range("N102").Value = -14803.8

$ws.Range("H132").Value = 7680.5806
$ws.Range("I132").Value = 3445.1365
$ws.Range("K132").Value = 10335.4095
$ws.Range("M132").Value = -7805.4095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10359.8125
$ws.Range("I7").Value = 4699.5
$ws.Range("J7").Value = 12246.583
$ws.Range("K7").Value = 4699.5
$ws.Range("L7").Value = 12246.583
$ws.Range("M7").Value = -4587.5
$ws.Range("N7").Value = -12470.583

$ws.Range("H16").Value = 1173.5366
$ws.Range("I16").Value = 1170.2413
$ws.Range("K16").Value = 1170.2413
$ws.Range("M16").Value = -1000.2413

$ws.Range("H126").Value = 10359.8125
$ws.Range("I126").Value = 4699.5
$ws.Range("J126").Value = 12246.583
$ws.Range("K126").Value = 14098.5
$ws.Range("L126").Value = 36739.749
$ws.Range("M126").Value = -11628.5
$ws.Range("N126").Value = -41679.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2109.3635
$ws.Range("J96").Value = 2572
$ws.Range("L96").Value = 2572
$ws.Range("N96").Value = -5318

$ws.Range("H100").Value = 791.3
$ws.Range("J100").Value = 800
$ws.Range("L100").Value = 1600
$ws.Range("N100").Value = -2682

$ws.Range("H115").Value = 355000
$ws.Range("J115").Value = 355000
$ws.Range("L115").Value = 355000
$ws.Range("N115").Value = -358134
